$wb = $excel.ActiveWorkbook

$wsServer = $wb.Worksheets.Item("Server")
$wsClient = $wb.Worksheets.Item("Client")

# ---------------------------------------------------------------------------
# Server sheet (sheet1)
# ---------------------------------------------------------------------------

# Copy the existing highlight format (column A marker cells) onto the rows
# that are newly marked in this revision.
$wsServer.Range("A2").Copy() | Out-Null
$wsServer.Range("A6").PasteSpecial(-4122) | Out-Null
$wsServer.Range("A9").PasteSpecial(-4122) | Out-Null
$wsServer.Range("A10").PasteSpecial(-4122) | Out-Null

# New row describing the "Desconexión de usuario" request code.
$wsServer.Range("B11").Value = "Desconexión de usuario"
$wsServer.Range("C11").Value = 110
$wsServer.Range("E11").Value = "UTF"
$wsServer.Range("F11").Value = "Global"

# ---------------------------------------------------------------------------
# Client sheet (sheet2)
# ---------------------------------------------------------------------------

$wsClient.Range("A2").Copy() | Out-Null
$wsClient.Range("A4").PasteSpecial(-4122) | Out-Null
$wsClient.Range("A5").PasteSpecial(-4122) | Out-Null
$wsClient.Range("A8").PasteSpecial(-4122) | Out-Null

# New row describing the "Usuario eliminado" response code.
$wsClient.Range("A12").PasteSpecial(-4122) | Out-Null
$wsClient.Range("B12").Value = "Usuario eliminado"
$wsClient.Range("C12").Value = 20
$wsClient.Range("E12").Value = "Sala"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selection / active cell bookkeeping (matches the recorded cursor position
# left behind by the author when they saved the workbook).
# ---------------------------------------------------------------------------

$wsServer.Activate()
$wsServer.Range("A9").Select() | Out-Null

$wsClient.Activate()
$wsClient.Range("A5").Select() | Out-Null
